# Insert a new row before row 307, shifting existing rows 307..386 down to
# 308..387, then populate the newly inserted row 307 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(307).Insert()

$ws.Cells.Item(307, 1).Value = 3
$ws.Cells.Item(307, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(307, 3).Value = "Coquimbo"
$ws.Cells.Item(307, 4).Value = 44543
$ws.Cells.Item(307, 5).Value = 5
$ws.Cells.Item(307, 6).Value = 100112045
$ws.Cells.Item(307, 7).Value = "Zapallo"
$ws.Cells.Item(307, 8).Value = "Camote"
$ws.Cells.Item(307, 9).Value = "1a nueva(o)"
$ws.Cells.Item(307, 10).Value = 120
$ws.Cells.Item(307, 11).Value = 700
$ws.Cells.Item(307, 12).Value = 700
$ws.Cells.Item(307, 13).Value = 700
$ws.Cells.Item(307, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(307, 15).Value = "Provincia de Talca"
$ws.Cells.Item(307, 16).Value = 700
$ws.Cells.Item(307, 17).Value = 1
$ws.Cells.Item(307, 18).Value = "Hortaliza"
